$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results for the 380 kV case (Case_2_91), row 2..25
# Columns B,C,E,F,G,J,L,M,O change; A,D,H,I,K,N stay at 0 (unchanged).
$data = @{
    2 = @{ "B"="1.225757106164735"; "C"="0.3129439817075621"; "E"="0.2308316458072017"; "F"="2.006044535286037"; "G"="0.002464079327741438"; "J"="0.06201190636205212"; "L"="0.4281784128119028"; "M"="0.341158130737746"; "O"="3.491593901693278" }
    3 = @{ "B"="1.125010755990672"; "C"="0.3034362689975438"; "E"="0.2321710535730546"; "F"="2.015264558025514"; "G"="0.002467035911406694"; "J"="0.06055202878836496"; "L"="0.4222000518085309"; "M"="0.3233105861790193"; "O"="3.522126831721948" }
    4 = @{ "B"="1.063273409146063"; "C"="0.2975807635467618"; "E"="0.2330570016367943"; "F"="2.022073380476684"; "G"="0.002468948076384555"; "J"="0.0596605893909441"; "L"="0.418654725366622"; "M"="0.3124264771974623"; "O"="3.543295530394403" }
    5 = @{ "B"="1.038146953549472"; "C"="0.2951903121779367"; "E"="0.2334340414005487"; "F"="2.02513662072154"; "G"="0.002469751716023572"; "J"="0.05929859681075555"; "L"="0.4172416817385738"; "M"="0.3080101121505976"; "O"="3.552530103587898" }
    6 = @{ "B"="1.033976703016265"; "C"="0.2947931258701573"; "E"="0.2334976163214639"; "F"="2.025662699365263"; "G"="0.00246988663672388"; "J"="0.0592385663707411"; "L"="0.4170089668854473"; "M"="0.3072779344095053"; "O"="3.554100207683689" }
    7 = @{ "B"="1.062934413344749"; "C"="0.2975485422029891"; "E"="0.2330620216686548"; "F"="2.022113523928191"; "G"="0.002468958815459978"; "J"="0.05965570221671612"; "L"="0.4186355399636881"; "M"="0.3123668392325243"; "O"="3.543417609631192" }
    8 = @{ "B"="1.190995453236155"; "C"="0.3096695029669974"; "E"="0.2312803023162502"; "F"="2.008985389032816"; "G"="0.002465078709076623"; "J"="0.06150753817704668"; "L"="0.4260911439647828"; "M"="0.3349890384593976"; "O"="3.501618673010455" }
    9 = @{ "B"="1.443032362522615"; "C"="0.3332915797961675"; "E"="0.228289255664075"; "F"="1.99234949350361"; "G"="0.002458234627843001"; "J"="0.06517672736841718"; "L"="0.4417002338209102"; "M"="0.3799303418044673"; "O"="3.43889709743766" }
    10 = @{ "B"="1.628706630127397"; "C"="0.3505502377410892"; "E"="0.2263965285943339"; "F"="1.985684814080102"; "G"="0.002453667819354884"; "J"="0.0678939348204679"; "L"="0.4537638965355342"; "M"="0.4132910310069136"; "O"="3.404597142864418" }
    11 = @{ "B"="1.713274262118148"; "C"="0.358379397514625"; "E"="0.2256012839887767"; "F"="1.983860927804642"; "G"="0.002451689477697038"; "J"="0.06913441650097241"; "L"="0.4593799274412476"; "M"="0.4285399715642058"; "O"="3.391562097635045" }
    12 = @{ "B"="1.74531149568071"; "C"="0.3613408053608111"; "E"="0.2253095732097794"; "F"="1.983344029084975"; "G"="0.002450954509044548"; "J"="0.06960475773543351"; "L"="0.4615248619620473"; "M"="0.4343245999883436"; "O"="3.386996207482383" }
    13 = @{ "B"="1.7384111372848"; "C"="0.3607031638748595"; "E"="0.2253719793085409"; "F"="1.983447622994547"; "G"="0.002451112167946569"; "J"="0.0695034353261974"; "L"="0.4610621022838757"; "M"="0.4330783286767428"; "O"="3.387963074666686" }
    14 = @{ "B"="1.715909726877157"; "C"="0.3586231019669981"; "E"="0.2255770959029402"; "F"="1.983814919272888"; "G"="0.002451628727320513"; "J"="0.06917309998538457"; "L"="0.4595560274027406"; "M"="0.4290156740647717"; "O"="3.391179034639663" }
    15 = @{ "B"="1.702128649775887"; "C"="0.3573485657497599"; "E"="0.2257039631757145"; "F"="1.98406253059737"; "G"="0.002451946980520666"; "J"="0.06897083681985805"; "L"="0.4586358867016429"; "M"="0.4265284980821704"; "O"="3.39319714036958" }
    16 = @{ "B"="1.623181883889572"; "C"="0.3500381280409783"; "E"="0.2264498208598855"; "F"="1.985828321708823"; "G"="0.002453799097597388"; "J"="0.06781295204621074"; "L"="0.4533994409869848"; "M"="0.4122959190853379"; "O"="3.405500768816864" }
    17 = @{ "B"="1.574775978887828"; "C"="0.345547677786584"; "E"="0.2269242061888388"; "F"="1.987221003293257"; "G"="0.002454960651072861"; "J"="0.06710373076090548"; "L"="0.4502197668862493"; "M"="0.4035831666275982"; "O"="3.413707057661497" }
    18 = @{ "B"="1.546944015393422"; "C"="0.3429628395940654"; "E"="0.2272032518811393"; "F"="1.988135726893816"; "G"="0.002455638079777584"; "J"="0.06669622271255804"; "L"="0.4484029803542313"; "M"="0.3985787109922114"; "O"="3.418668773307559" }
    19 = @{ "B"="1.537522327737634"; "C"="0.3420873111509195"; "E"="0.2272987962002571"; "F"="1.988464960759842"; "G"="0.002455869050591949"; "J"="0.06655832027586683"; "L"="0.447789927270037"; "M"="0.396885481040222"; "O"="3.420390205555606" }
    20 = @{ "B"="1.579927863897979"; "C"="0.3460259073384861"; "E"="0.226873066421966"; "F"="1.987060982949643"; "G"="0.002454836036168904"; "J"="0.06717918565618675"; "L"="0.4505569998461283"; "M"="0.4045099432632853"; "O"="3.412808465189698" }
    21 = @{ "B"="1.722518587879733"; "C"="0.3592341582117058"; "E"="0.225516592436847"; "F"="1.983702318890877"; "G"="0.00245147661686556"; "J"="0.06927011154813556"; "L"="0.4599979037399322"; "M"="0.4302087002373582"; "O"="3.390224375012821" }
    22 = @{ "B"="1.815786710852763"; "C"="0.3678470522160353"; "E"="0.2246850179024609"; "F"="1.982520119646964"; "G"="0.002449363701625992"; "J"="0.07064011998817676"; "L"="0.4662744668975165"; "M"="0.4470635503233353"; "O"="3.377622582080392" }
    23 = @{ "B"="1.766001243213395"; "C"="0.363252027265105"; "E"="0.225123824645884"; "F"="1.983058379576747"; "G"="0.002450483863034183"; "J"="0.06990861535807369"; "L"="0.4629148688494524"; "M"="0.4380624824919934"; "O"="3.384150623248672" }
    24 = @{ "B"="1.577598705177422"; "C"="0.345809709771487"; "E"="0.2268961670374914"; "F"="1.987132972883927"; "G"="0.002454892344633954"; "J"="0.06714507176990736"; "L"="0.4504045017697393"; "M"="0.4040909331269873"; "O"="3.413213959444846" }
    25 = @{ "B"="1.374757731802788"; "C"="0.3269176426640001"; "E"="0.2290447592377216"; "F"="1.995874310403977"; "G"="0.00246000474448589"; "J"="0.06418023623729852"; "L"="0.4373724507929779"; "M"="0.367711644325432"; "O"="3.453799484884655" }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = [double]$data[$row][$col]
    }
}
